# Practiced IELTS Ready MockTest 3 Reading.
# Fill in rows 41 and 42 of Sheet1 (Table1) with the new mock-test
# practice entries, and update the sheet view's selection/scroll
# position to reflect where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed row 41/42 formatting (styles, number formats, borders) from row 40,
# the prior fully-filled-in row, then overwrite with the new row's data.
$ws.Range("C40:L40").Copy()
$ws.Range("C41:L41").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C40:L40").Copy()
$ws.Range("C42:L42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 41: IELTS_Ready_MockTest_2 - Listening 29 -> 6.5, Reading 26 -> 6
$ws.Range("C41").Value = 38
$ws.Range("D41").Value = 45522
$ws.Range("E41").Value = "IELTS_Ready_MockTest_2"
$ws.Range("F41").Value = 29
$ws.Range("G41").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"
$ws.Range("H41").Value = 26
$ws.Range("I41").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"

# Row 42: IELTS_Ready_MockTest_3 - Reading only, 31 -> 7
$ws.Range("C42").Value = 39
$ws.Range("D42").Value = 45525
$ws.Range("E42").Value = "IELTS_Ready_MockTest_3"
$ws.Range("H42").Value = 31
$ws.Range("I42").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"

# Update the view: scrolled down a bit further, and selection moved to H43
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("H43").Select()
